# "custom accuracy + data-1000" edit:
#  - Round the numeric readings in row 5 (B5:AH5) to 2 decimal places
#    (custom accuracy: the raw simulation export used 3 decimals).
#  - Drop the now-redundant row 6 (trimming the sample set), which also
#    shrinks the used range down to A1:AH5.
#  - Two columns (B and W) had their single widest entry in row 5; once
#    that entry is rounded down to 2 decimals (and row 6, its other wide
#    neighbour, is gone) their auto-fit width narrows by one unit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the numeric values in row 5, columns B..AH, to 2 decimal places
# (round-half-away-from-zero, matching Excel's ROUND()).
$rowRange = $ws.Range("B5:AH5")
foreach ($cell in $rowRange.Cells) {
    $val = $cell.Value2
    if ($null -ne $val) {
        $scaled = $val * 100
        if ($scaled -ge 0) {
            $rounded = [Math]::Floor($scaled + 0.5) / 100
        } else {
            $rounded = [Math]::Ceiling($scaled - 0.5) / 100
        }
        $cell.Value2 = $rounded
    }
}

# Delete row 6 entirely (rows below shift up; none here).
$ws.Rows.Item(6).Delete()

# Narrow columns B and W by one character, mirroring the auto-fit
# shrink once the row-5/row-6 wide values are gone. ColumnWidth is
# stored internally with a constant ~0.8333 padding added back on
# export, so back that off to land exactly on the target width.
$padding = 0.8333333333333334
$ws.Columns.Item(2).ColumnWidth = 7 - $padding
$ws.Columns.Item(23).ColumnWidth = 7 - $padding
